$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.168994188308716
$ws.Range("B1").Value = 2.139285326004028
$ws.Range("C1").Value = 3.518226146697998
$ws.Range("D1").Value = 3.442429780960083
$ws.Range("E1").Value = 1.169265985488892
